$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the candidate record in row 2 with new generated values
$ws.Range("A2").Value = "FCLXw697"
$ws.Range("B2").Value = 23081037
$ws.Range("C2").Value = "hwgktnx59"
$ws.Range("D2").Value = "h8!7#QBp"
$ws.Range("F2").Value = "RpFFQEZI"
$ws.Range("G2").Value = "Feuq"
